$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (PEDC), shifting the old
# TOTAL_BCQ column (and its values) from E to F.
$ws.Columns.Item(5).Insert()

# New header for the inserted column
$ws.Range("E1").Value = "PEDC"

# Row data: Hour, SCPC(B), KSPC(C), EDC(D), PEDC(E, new), TOTAL_BCQ(F)
$data = @(
    @(2, 16500, 20000, 20000, 10000, 66500),
    @(3, 13500, 20000, 20000, 10000, 63500),
    @(4, 12500, 10000, 20000, 0, 42500),
    @(5, 12500, 10000, 20000, 0, 42500),
    @(6, 12500, 10000, 0, 0, 22500),
    @(7, 12500, 10000, 0, 0, 22500),
    @(8, 12500, 10000, 0, 0, 22500),
    @(9, 12500, 10000, 0, 0, 22500),
    @(10, 12500, 10000, 0, 0, 22500),
    @(11, 12500, 10000, 10000, 0, 32500),
    @(12, 12500, 10000, 10000, 0, 32500),
    @(13, 12500, 10000, 10000, 0, 32500),
    @(14, 12500, 10000, 10000, 0, 32500),
    @(15, 12500, 10000, 20000, 10000, 52500),
    @(16, 25000, 20000, 20000, 10000, 75000),
    @(17, 25000, 20000, 20000, 10000, 75000),
    @(18, 25000, 20000, 20000, 10000, 75000),
    @(19, 25000, 20000, 20000, 10000, 75000),
    @(20, 25000, 20000, 20000, 10000, 75000),
    @(21, 25000, 20000, 20000, 10000, 75000),
    @(22, 25000, 20000, 20000, 10000, 75000),
    @(23, 25000, 20000, 20000, 10000, 75000),
    @(24, 22000, 20000, 20000, 10000, 72000),
    @(25, 12500, 10000, 20000, 10000, 52500)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
